$d = $word.ActiveDocument

# Locate "IDG PA28X" without mutating anything yet, so we know exactly
# where the trailing "X" (to be dropped) and the run boundary sit.
$findRng = $d.Content.Duplicate
$null = $findRng.Find.Execute("IDG PA28X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleStart = $findRng.Start
$titleEnd = $findRng.End
$xStart = $titleEnd - 1
$xEnd = $titleEnd

# The document carries a stray "_GoBack" bookmark left over from the
# previous edit session (right after "Airspeed Indicator"). Remove it so
# it can be re-added at the new edit location, mirroring what Word does.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Drop the trailing "X": "IDG PA28X" -> "IDG PA28".
$xRng = $d.Range($xStart, $xEnd)
$xRng.Text = ""

# Re-drop the "_GoBack" bookmark right after the edited run, i.e. right
# before the ": Cockpit" run that follows it.
$bmRng = $d.Range($xStart, $xStart)
$d.Bookmarks.Add("_GoBack", $bmRng)
